$d = $word.ActiveDocument

# Update the date heading.
$d.Paragraphs.Item(1).Range.Text = "2024-09-13 Friday"

# Update the worksheet answers. Values are written by (row, column)
# table coordinates rather than a global find/replace, because some
# new values coincide with other cells' old values (e.g. "99÷4=24, 3"
# is both the new answer for cell (1,1) and the old answer for cell
# (5,3)) -- a sequential text replace would cascade incorrectly.
$t = $d.Tables.Item(1)

$answers = @{
    1 = @("99÷4=24, 3", "54÷8=6, 6", "40÷6=6, 4", "59÷8=7, 3", "85÷6=14, 1")
    5 = @("91÷2=45, 1", "74÷3=24, 2", "30÷7=4, 2", "95÷9=10, 5", "93÷2=46, 1")
    9 = @("53÷4=13, 1", "69÷3=23, 0", "84÷5=16, 4", "95÷7=13, 4", "68÷3=22, 2")
    13 = @("30÷3=10, 0", "56÷8=7, 0", "43÷6=7, 1", "77÷4=19, 1", "36÷2=18, 0")
    17 = @("48÷3=16, 0", "68÷2=34, 0", "51÷9=5, 6", "91÷3=30, 1", "11÷3=3, 2")
}

foreach ($row in $answers.Keys) {
    $cols = $answers[$row]
    for ($c = 0; $c -lt $cols.Length; $c++) {
        $t.Cell($row, $c + 1).Range.Text = $cols[$c]
    }
}
